# Apply "enter vpn app connection page execution" edits to sprint_38.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Day 4 block: Total testcase Written bumped from 538 to 539
$ws.Range("C17").Value = 539

# Day 5 block (rows 24-27): fill in execution numbers that were previously empty
$ws.Range("C25").Value = 539
$ws.Range("C26").Value = 802
$ws.Range("C27").Value = 511

# Update the view so it reflects where the user was working (Day 5 block)
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C27").Select()
